$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample row ("Shrouk" / "pasdasdoj") with a real record:
# A2 becomes a mailto hyperlink displaying the email address, B2 becomes a phone/ID number.
$ws.Range("A2").Value = "shrouk.ali@gmail.com"
$ws.Range("B2").Value = 521993
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:shrouk.ali@gmail.com")

# Move the active cell selection to B2
$ws.Range("B2").Select() | Out-Null
